$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: style change (s=3 -> s=2), reorder A7, fill M7 with new Fix text ---
$ws.Range("A6:Q6").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$ws.Range("A7").Value = 'AU-5 a,AU-5 b'
$ws.Range("M7").Value = 'Configure Red Hat Enterprise Linux 9 to shutdown when auditing failures occur.
If the auditd daemon is configured to use the augenrules program to read
audit rules during daemon startup (the default), add the following line to
the bottom of "/etc/audit/rules.d/immutable.rules":
-f 2
If the auditd daemon is configured to use the auditctl utility to read
audit rules during daemon startup, add the following line to the
bottom of the /etc/audit/audit.rules file:
-f 2'

# --- Row 12: K12 wording fix ---
$ws.Range("K12").Value = 'For every temporary account, run the following command
to obtain its account aging and expiration information:
 $ sudo chage -l  USER  
Verify each of these accounts has an expiration date set as documented.

If any temporary accounts have no expiration date set or do not expire within a documented time frame, then this is a finding.'

# --- Column A reorder-only rows ---
$ws.Range("A13").Value = 'MA-4 (1) (a),AU-12 a,AU-7 a,CM-6 b,AU-7 (1),AU-6 (4),AU-3,AU-3 (1),AU-14 (1),CM-5 (1)'
$ws.Range("A37").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-14 (1),AU-3,AU-3 (1)'
$ws.Range("A38").Value = 'AU-4,AU-14 (1)'
$ws.Range("A44").Value = 'CM-6 b,AU-4 (1),AU-6 (4)'
$ws.Range("A48").Value = 'IA-2 (12),IA-2 (11)'
$ws.Range("A49").Value = 'IA-2 (1),IA-2 (12),IA-2 (11)'
$ws.Range("A50").Value = 'SI-6 d,SI-6 b,CM-3 (5)'
$ws.Range("A51").Value = 'SI-6 d,CM-3 (5)'
$ws.Range("A52").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A53").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A54").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A55").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A56").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A57").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A58").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A59").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A60").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A61").Value = 'AU-7 a,AU-12 c,AU-12 (3),AU-8 b,CM-6 b,AU-7 b,AU-12 a,CM-5 (1)'
$ws.Range("A62").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A63").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A64").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A65").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A66").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A67").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A68").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A69").Value = 'SI-6 a,CM-3 (5)'
$ws.Range("A82").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A83").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A84").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A85").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A92").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A93").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A94").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A95").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A96").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A97").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A98").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A99").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A100").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A101").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A102").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A103").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A104").Value = 'AC-11 b,AC-11 (1)'
$ws.Range("A106").Value = 'AC-11 b,AC-6 (10)'
$ws.Range("A107").Value = 'AC-11 b,AC-11 a'
$ws.Range("A108").Value = 'AC-11 b,AC-11 a'
$ws.Range("A109").Value = 'AC-11 b,AC-11 a'
$ws.Range("A110").Value = 'AC-11 b,AC-11 a'
$ws.Range("A113").Value = 'AC-17 (2),SC-13,MA-4 c,SC-8'
$ws.Range("A114").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A115").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A116").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A117").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A118").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A119").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A120").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A121").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A122").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A123").Value = 'SC-28,SC-28 (1)'
$ws.Range("A133").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A134").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A135").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A136").Value = 'AC-6 (10),AC-3 (4)'
$ws.Range("A137").Value = 'AC-6 (10),AC-3 (4)'
$ws.Range("A139").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A140").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A141").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A142").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A146").Value = 'IA-5 (1) (c),IA-7'
$ws.Range("A150").Value = 'AC-11 a,AC-11 (1)'
$ws.Range("A152").Value = 'IA-7,CM-6 b'
$ws.Range("A153").Value = 'IA-7,CM-6 b'
$ws.Range("A154").Value = 'IA-7,CM-6 b'
$ws.Range("A157").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A158").Value = 'AU-3,AU-12 a,AU-12 c,MA-4 (1) (a)'
$ws.Range("A159").Value = 'AU-3,AU-12 a,AU-12 c,MA-4 (1) (a)'
$ws.Range("A160").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A161").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A162").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A163").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A164").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A165").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A166").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A167").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A168").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A169").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A170").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A171").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A172").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A173").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A174").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A175").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A176").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A177").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A178").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A179").Value = 'AU-3,AU-3 (1),MA-4 (1) (a)'
$ws.Range("A180").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A181").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A182").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A183").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A184").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A185").Value = 'MA-4 (1) (a),AU-3 (1),AU-12 c'
$ws.Range("A186").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A187").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A188").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A189").Value = 'AU-3,AU-3 (1),AU-12 c,MA-4 (1) (a)'
$ws.Range("A190").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A191").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A192").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A193").Value = 'MA-4 (1) (a),AU-12 a,AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A194").Value = 'MA-4 (1) (a),AU-12 a,AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A195").Value = 'MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A196").Value = 'IA-2 (4),IA-2 (3),IA-2 (1),IA-2 (2)'
$ws.Range("A197").Value = 'IA-2 (5),IA-2 (2),IA-2 (3),IA-2,IA-2 (4)'
$ws.Range("A198").Value = 'IA-2 (5),IA-2 (2),IA-2 (3),IA-2,IA-2 (4)'
$ws.Range("A201").Value = 'IA-11,AC-3 (4)'
$ws.Range("A206").Value = 'SC-8 (2),SC-8 (1),SC-8'
$ws.Range("A207").Value = 'SC-8 (2),SC-8 (1),SC-8'
$ws.Range("A208").Value = 'AC-18 (1),SC-8 (1),SC-8'
$ws.Range("A209").Value = 'AC-11 a,AC-11 (1)'
$ws.Range("A210").Value = 'AC-11 a,AC-11 (1)'
$ws.Range("A213").Value = 'IA-5 (1) (c),CM-7 a,CM-6 b'
$ws.Range("A215").Value = 'AU-12 a,CM-6 b'
$ws.Range("A216").Value = 'SC-5,CM-6 b,SC-5 (2)'
$ws.Range("A220").Value = 'IA-8,AU-3 (1),IA-2'
$ws.Range("A232").Value = 'SC-2,SI-16,CM-6 b'
$ws.Range("A243").Value = 'CM-6 b,IA-5 (1) (a),IA-5 (1) (b)'
$ws.Range("A250").Value = 'SC-4,CM-6 b'
$ws.Range("A254").Value = 'IA-2 (12),IA-2 (11)'
$ws.Range("A260").Value = 'IA-2 (5),CM-6 b'
$ws.Range("A267").Value = 'CM-5 (3),CM-6 b'
$ws.Range("A275").Value = 'CM-7 a,CM-6 b'
$ws.Range("A277").Value = 'CM-7 b,CM-7 a'
$ws.Range("A278").Value = 'CM-7 b,CM-7 a'
$ws.Range("A292").Value = 'AC-6 (9),AC-2 (4),AU-12 c'
$ws.Range("A339").Value = 'IA-5 (1) (c),CM-6 b'
$ws.Range("A374").Value = 'AU-3,CM-6 b'
$ws.Range("A377").Value = 'SC-3,CM-6 b'
$ws.Range("A382").Value = 'CM-6 b,CM-7 b,AC-17 (1),AC-17 (9)'
$ws.Range("A383").Value = 'CM-6 b,CM-7 b,AC-17 (1)'
$ws.Range("A412").Value = 'IA-3,CM-6 b'
$ws.Range("A413").Value = 'IA-3,CM-6 b'
$ws.Range("A422").Value = 'SC-3,CM-6 b'
$ws.Range("A424").Value = 'SC-2,CM-6 b'
$ws.Range("A425").Value = 'SC-2,CM-6 b'
$ws.Range("A429").Value = 'SC-3,CM-6 b'
$ws.Range("A439").Value = 'IA-3,CM-6 b'
$ws.Range("A440").Value = 'IA-3,CM-6 b'
$ws.Range("A458").Value = 'MA-4 e,MA-4 (7),AC-12,SC-10'
$ws.Range("A459").Value = 'AC-12,SC-10'
$ws.Range("A460").Value = 'AC-12,SC-10'
$ws.Range("A462").Value = 'AC-17 (2),SC-8 (1),SC-8'
$ws.Range("A483").Value = 'CM-7 b,AC-17 (1)'
$ws.Range("A497").Value = 'AU-4 (1),AU-4'

# --- Row 124: H124 fill with Vul Discussion text ---
$ws.Range("H124").Value = 'A replay attack may enable an unauthorized user to gain access to Red Hat Enterprise Linux 9. Authentication sessions between the authenticator and Red Hat Enterprise Linux 9 validating the user credentials must not be vulnerable to a replay attack.
An authentication process resists replay attacks if it is impractical to achieve a successful authentication by recording and replaying a previous authentication message.
A privileged account is any information system account with authorizations of a privileged user.
Techniques used to address this include protocols using nonces (e.g., numbers generated for a specific one-time use) or challenges (e.g., TLS, WS_Security). Additional techniques include time-synchronous or challenge-response one-time authenticators.'

# --- Row 293: style change (s=3 -> s=2), K293 wording fix, fill M293 ---
$ws.Range("A6:Q6").Copy()
$ws.Range("A293:Q293").PasteSpecial(-4122)
$ws.Range("K293").Value = 'Find the list of alias maps used by the Postfix mail server:
 $ sudo postconf alias_maps 
Query the Postfix alias maps for an alias for the  root  user:
 $ sudo postmap -q root hash:/etc/aliases 
The output should return an alias.

If the alias is not set, then this is a finding.'
$ws.Range("M293").Value = 'Configure a valid email address as an alias for the root account.
Append the following line to "/etc/aliases":
root: system.administrator@mail.mil
Then, run the following command:
$ sudo newaliases'

# --- Row 294: style change (s=3 -> s=2), reorder A294, fill M294 ---
$ws.Range("A6:Q6").Copy()
$ws.Range("A294:Q294").PasteSpecial(-4122)
$ws.Range("A294").Value = 'AU-5 a,AU-5 (1)'
$ws.Range("M294").Value = 'Configure "auditd" service to notify the SA and ISSO in the event of an audit processing failure.
Edit the following line in "/etc/audit/auditd.conf" to ensure that administrators are notified via email for those situations:
action_mail_acct = root'

